$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 184 (shifts existing rows 184:211 down to 185:212)
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new record
$ws.Cells.Item(184, 1).Value = 5
$ws.Cells.Item(184, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(184, 3).Value = "Maule"
$ws.Cells.Item(184, 4).Value = 44776
$ws.Cells.Item(184, 5).Value = 7
$ws.Cells.Item(184, 6).Value = 100112017
$ws.Cells.Item(184, 7).Value = "Apio"
$ws.Cells.Item(184, 8).Value = "Americana (o)"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 500
$ws.Cells.Item(184, 11).Value = 9000
$ws.Cells.Item(184, 12).Value = 9000
$ws.Cells.Item(184, 13).Value = 9000
$ws.Cells.Item(184, 14).Value = "$/docena de matas"
$ws.Cells.Item(184, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(184, 16).Value = 1500
$ws.Cells.Item(184, 17).Value = 6
$ws.Cells.Item(184, 18).Value = "Hortaliza"
